# NATMI re-run following Dr Hou advice: adds two new target clusters (M2 used
# also as a target, plus a brand-new "sCs" cluster) and refreshes every
# expression / specificity statistic for the Ebi3->Il27ra edge table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Ebi3", "Il27ra", "ECs", 1, 0.3333333333333333, 0.1351216666666667, 0.405365, 0.01567245347423418, 0.01567245347423418, 3, 1, 1.423311333333333, 4.269934, 0.505384319942809, 0.505384319942809, 0.1923201995455556, 1.73088179591, 0.007920612240911154, 0.007920612240911154),
    @("ECs", "Ebi3", "Il27ra", "FAPs", 1, 0.3333333333333333, 0.1351216666666667, 0.405365, 0.01567245347423418, 0.01567245347423418, 3, 1, 1.169098333333333, 3.507295, 0.4151192731348574, 0.4151192731348574, 0.1579705152972222, 1.421734637675, 0.006505937494463964, 0.006505937494463964),
    @("ECs", "Ebi3", "Il27ra", "M2", 1, 0.3333333333333333, 0.1351216666666667, 0.405365, 0.01567245347423418, 0.01567245347423418, 2, 0.6666666666666666, 0.192572, 0.577716, 0.06837778002659521, 0.06837778002659523, 0.02602064959333333, 0.23418584634, 0.001071647576138233, 0.001071647576138233),
    @("ECs", "Ebi3", "Il27ra", "sCs", 1, 0.3333333333333333, 0.1351216666666667, 0.405365, 0.01567245347423418, 0.01567245347423418, 1, 0.3333333333333333, 0.03131333333333333, 0.09394, 0.01111862689573831, 0.01111862689573831, 0.004231109788888889, 0.0380799881, 0.0001742561627208275, 0.0001742561627208275),
    @("M2", "Ebi3", "Il27ra", "ECs", 3, 1, 8.486481, 25.459443, 0.9843275465257658, 0.9843275465257658, 3, 1, 1.423311333333333, 4.269934, 0.505384319942809, 0.505384319942809, 12.078904587418, 108.710141286762, 0.4974637077018978, 0.4974637077018978),
    @("M2", "Ebi3", "Il27ra", "FAPs", 3, 1, 8.486481, 25.459443, 0.9843275465257658, 0.9843275465257658, 3, 1, 1.169098333333333, 3.507295, 0.4151192731348574, 0.4151192731348574, 9.921530792965, 89.29377713668501, 0.4086133356403934, 0.4086133356403934),
    @("M2", "Ebi3", "Il27ra", "M2", 3, 1, 8.486481, 25.459443, 0.9843275465257658, 0.9843275465257658, 2, 0.6666666666666666, 0.192572, 0.577716, 0.06837778002659521, 0.06837778002659523, 1.634258619132, 14.708327572188, 0.06730613245045698, 0.067306132450457),
    @("M2", "Ebi3", "Il27ra", "sCs", 3, 1, 8.486481, 25.459443, 0.9843275465257658, 0.9843275465257658, 1, 0.3333333333333333, 0.03131333333333333, 0.09394, 0.01111862689573831, 0.01111862689573831, 0.2657400083799999, 2.39166007542, 0.01094437073301749, 0.01094437073301749),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

Write-Output "Updated $($data.Count) data rows (rows $startRow..$($startRow + $data.Count - 1))."
